$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values in column D are written as literal text
# (matching the source data which stores prices as plain strings, e.g. "70.306.87"),
# rather than being auto-converted to floating point numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '70.306.87'
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('D3').Value = '3.550.69'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '617.03'
$ws.Range('E5').Value = '  +5.42%  '
$ws.Range('D6').Value = '189.03'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').Value = '0.640'
$ws.Range('E7').Value = '  +2.92%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').Value = '0.666'
$ws.Range('E10').Value = '  +1.67%  '
$ws.Range('D11').Value = '54.05'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  -3.75%  '
$ws.Range('D13').Value = '9.77'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').Value = '4.116.47'
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').Value = '614.50'
$ws.Range('E15').Value = '  +8.20%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = '12.88'
$ws.Range('E16').Value = '  +3.81%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '70.360.98'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('D18').Value = '19.18'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('D19').Value = '3.554.12'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').Value = '17.75'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').Value = '105.90'
$ws.Range('E23').Value = '  +11.42%  '
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').Value = '5.15'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').Value = '3.05'
$ws.Range('E26').Value = '  +3.94%  '
$ws.Range('D27').Value = '11.02'
$ws.Range('E27').Value = '  -3.24%  '
$ws.Range('D28').Value = '10.16'
$ws.Range('E28').Value = '  +10.81%  '
$ws.Range('D29').Value = '34.63'
$ws.Range('E29').Value = '  +6.76%  '
$ws.Range('D30').Value = '7.14'
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').Value = '12.59'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').Value = '0.118'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('D33').Value = '64.38'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').Value = '3.75'
$ws.Range('E34').Value = '  +14.60%  '
$ws.Range('D35').Value = '3.18'
$ws.Range('E35').Value = '  -5.60%  '
$ws.Range('D36').Value = '539.67'
$ws.Range('E36').Value = '  -2.70%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '0.403'
$ws.Range('E38').Value = '  -3.46%  '
$ws.Range('D39').Value = '37.42'
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('D40').Value = '0.0₃0785'
$ws.Range('E40').Value = '  -3.21%  '
$ws.Range('D41').Value = '3.58'
$ws.Range('E41').Value = '  +4.14%  '
$ws.Range('D42').Value = '3.544.16'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').Value = '0.140'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('E44').Value = '  +4.81%  '
$ws.Range('D45').Value = '2.98'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').Value = '0.145'
$ws.Range('E46').Value = '  +5.03%  '
$ws.Range('D47').Value = '3.37'
$ws.Range('E47').Value = '  -3.25%  '
$ws.Range('D48').Value = '9.04'
$ws.Range('E48').Value = '  -4.52%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '133.77'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('E51').Value = '  -4.92%  '
